$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 8984.286
$ws.Range("I6").Value = 9143.549999999999
$ws.Range("J6").Value = 5799
$ws.Range("K6").Value = 27430.65
$ws.Range("L6").Value = 17397
$ws.Range("M6").Value = -27318.65
$ws.Range("N6").Value = -17621
$ws.Range("H28").Value = 771.5454999999999
$ws.Range("I28").Value = 398.55554
$ws.Range("K28").Value = 398.55554
$ws.Range("M28").Value = 86.44445999999999
$ws.Range("H111").Value = 1428.8
$ws.Range("I111").Value = 1411
$ws.Range("K111").Value = 4233
$ws.Range("M111").Value = -1166
$ws.Range("H113").Value = 6754.5713
$ws.Range("I113").Value = 7217.273
$ws.Range("K113").Value = 7217.273
$ws.Range("M113").Value = -3963.273
$ws.Range("H135").Value = 2540.7273
$ws.Range("I135").Value = 2494.4443
$ws.Range("J135").Value = 2749
$ws.Range("K135").Value = 22449.9987
$ws.Range("L135").Value = 24741
$ws.Range("M135").Value = -19914.9987
$ws.Range("N135").Value = -29811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5453
$ws.Range("I2").Value = 3922.6191
$ws.Range("K2").Value = 3922.6191
$ws.Range("M2").Value = -3809.6191
$ws.Range("H61").Value = 4665.5415
$ws.Range("I61").Value = 2735.7896
$ws.Range("K61").Value = 2735.7896
$ws.Range("M61").Value = -2523.7896
$ws.Range("H116").Value = 5453
$ws.Range("I116").Value = 3922.6191
$ws.Range("K116").Value = 3922.6191
$ws.Range("M116").Value = -1628.6191
$ws.Range("H122").Value = 2512.5454
$ws.Range("I122").Value = 2409.7646
$ws.Range("K122").Value = 7229.293799999999
$ws.Range("M122").Value = -4779.293799999999
$ws.Range("H136").Value = 4665.5415
$ws.Range("I136").Value = 2735.7896
$ws.Range("K136").Value = 8207.3688
$ws.Range("M136").Value = -5657.3688

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5453
$ws.Range("I3").Value = 3922.6191
$ws.Range("K3").Value = 3922.6191
$ws.Range("M3").Value = -3808.6191
$ws.Range("H58").Value = 27249.5
$ws.Range("J58").Value = 22999.666
$ws.Range("L58").Value = 22999.666
$ws.Range("N58").Value = -23587.666
$ws.Range("H74").Value = 59948.5
$ws.Range("J74").Value = 59948.5
$ws.Range("L74").Value = 59948.5
$ws.Range("N74").Value = -61820.5
$ws.Range("H77").Value = 59948.5
$ws.Range("J77").Value = 59948.5
$ws.Range("L77").Value = 179845.5
$ws.Range("N77").Value = -189205.5
$ws.Range("H94").Value = 683
$ws.Range("I94").Value = 201.3
$ws.Range("K94").Value = 201.3
$ws.Range("M94").Value = 249.7
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H139").Value = 59983.855
$ws.Range("J139").Value = 59983.855
$ws.Range("L139").Value = 59983.855
$ws.Range("N139").Value = -70263.85500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9605.609
$ws.Range("I99").Value = 6988.125
$ws.Range("K99").Value = 6988.125
$ws.Range("M99").Value = -5490.125
$ws.Range("H126").Value = 9605.609
$ws.Range("I126").Value = 6988.125
$ws.Range("K126").Value = 20964.375
$ws.Range("M126").Value = -18494.375
$ws.Range("H131").Value = 1723333.4
$ws.Range("J131").Value = 1723333.4
$ws.Range("L131").Value = 1723333.4
$ws.Range("N131").Value = -1733413.4
$ws.Range("H141").Value = 143298.5
$ws.Range("J141").Value = 245449.25
$ws.Range("L141").Value = 245449.25
$ws.Range("N141").Value = -255809.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 311.9375
$ws.Range("I2").Value = 59.77778
$ws.Range("K2").Value = 358.66668
$ws.Range("M2").Value = -245.66668
$ws.Range("H23").Value = 85.375
$ws.Range("J23").Value = 87.14286
$ws.Range("L23").Value = 261.42858
$ws.Range("N23").Value = -731.42858
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 30084
$ws.Range("I62").Value = 30084
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 30084
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -29398
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 30084
$ws.Range("I65").Value = 30084
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 90252
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -86820
$ws.Range("N65").ClearContents()
$ws.Range("H102").Value = 1693.1428
$ws.Range("I102").Value = 1595.6
$ws.Range("K102").Value = 1595.6
$ws.Range("M102").Value = 26.40000000000009
$ws.Range("H107").Value = 571.3077
$ws.Range("I107").Value = 420.85715
$ws.Range("K107").Value = 420.85715
$ws.Range("M107").Value = 1499.14285
$ws.Range("H113").Value = 288270.44
$ws.Range("I113").Value = 335824.16
$ws.Range("K113").Value = 335824.16
$ws.Range("M113").Value = -333654.16
$ws.Range("H121").Value = 83093
$ws.Range("J121").Value = 83093
$ws.Range("L121").Value = 83093
$ws.Range("N121").Value = -86587
$ws.Range("H133").Value = 100498
$ws.Range("J133").Value = 100498
$ws.Range("L133").Value = 100498
$ws.Range("N133").Value = -110618

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 832
$ws.Range("I7").Value = 832
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 832
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -720
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 2762.25
$ws.Range("I22").Value = 2180.4546
$ws.Range("K22").Value = 2180.4546
$ws.Range("M22").Value = -1885.4546
$ws.Range("H23").Value = 100000
$ws.Range("I23").Value = 100000
$ws.Range("K23").Value = 100000
$ws.Range("M23").Value = -99770
$ws.Range("H27").Value = 2762.25
$ws.Range("I27").Value = 2180.4546
$ws.Range("K27").Value = 2180.4546
$ws.Range("M27").Value = -2073.4546
$ws.Range("H46").Value = 1064.1428
$ws.Range("I46").Value = 1100
$ws.Range("K46").Value = 1100
$ws.Range("M46").Value = -912
$ws.Range("H61").Value = 2812.2222
$ws.Range("I61").Value = 2817.625
$ws.Range("K61").Value = 2817.625
$ws.Range("M61").Value = -2615.625
$ws.Range("H68").Value = 2757.9546
$ws.Range("I68").Value = 2583.9
$ws.Range("J68").Value = 4498.5
$ws.Range("K68").Value = 2583.9
$ws.Range("L68").Value = 4498.5
$ws.Range("M68").Value = -1834.9
$ws.Range("N68").Value = -5996.5
$ws.Range("H71").Value = 2757.9546
$ws.Range("I71").Value = 2583.9
$ws.Range("J71").Value = 4498.5
$ws.Range("K71").Value = 12919.5
$ws.Range("L71").Value = 22492.5
$ws.Range("M71").Value = -9175.5
$ws.Range("N71").Value = -29980.5
$ws.Range("H93").Value = 2687.3076
$ws.Range("I93").Value = 2837.0527
$ws.Range("J93").Value = 2280.8572
$ws.Range("K93").Value = 2837.0527
$ws.Range("L93").Value = 2280.8572
$ws.Range("M93").Value = -1589.0527
$ws.Range("N93").Value = -4776.8572
$ws.Range("H108").Value = 60626
$ws.Range("J108").Value = 60626
$ws.Range("L108").Value = 60626
$ws.Range("N108").Value = -68306
$ws.Range("H113").Value = 2812.2222
$ws.Range("I113").Value = 2817.625
$ws.Range("K113").Value = 2817.625
$ws.Range("M113").Value = -647.625
$ws.Range("H126").Value = 832
$ws.Range("I126").Value = 832
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2496
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -26
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 5939
$ws.Range("I136").Value = 4548.75
$ws.Range("K136").Value = 13646.25
$ws.Range("M136").Value = -11096.25
$ws.Range("H140").Value = 55969
$ws.Range("J140").Value = 61213.75
$ws.Range("L140").Value = 61213.75
$ws.Range("N140").Value = -71573.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2314.8333
$ws.Range("I81").Value = 926
$ws.Range("J81").Value = 2777.7778
$ws.Range("K81").Value = 1852
$ws.Range("L81").Value = 5555.5556
$ws.Range("M81").Value = -791
$ws.Range("N81").Value = -7677.5556
$ws.Range("H84").Value = 2314.8333
$ws.Range("I84").Value = 926
$ws.Range("J84").Value = 2777.7778
$ws.Range("K84").Value = 9260
$ws.Range("L84").Value = 27777.778
$ws.Range("M84").Value = -3956
$ws.Range("N84").Value = -38385.778
$ws.Range("H113").Value = 462
$ws.Range("J113").Value = 416.66666
$ws.Range("L113").Value = 1249.99998
$ws.Range("N113").Value = -5589.999980000001
$ws.Range("H126").Value = 5786.613
$ws.Range("I126").Value = 5666.8887
$ws.Range("J126").Value = 6594.75
$ws.Range("K126").Value = 17000.6661
$ws.Range("L126").Value = 19784.25
$ws.Range("M126").Value = -14530.6661
$ws.Range("N126").Value = -24724.25
$ws.Range("H136").Value = 2225.4814
$ws.Range("I136").Value = 1095.0454
$ws.Range("K136").Value = 3285.1362
$ws.Range("M136").Value = -735.1361999999999
$ws.Range("H139").Value = 46381.668
$ws.Range("J139").Value = 46381.668
$ws.Range("L139").Value = 46381.668
$ws.Range("N139").Value = -56661.668
